# LocPage & Add LocationPage Partially completed
# Adds two new worksheets (DeleteLocations, AddLocation) to the HRM test-data
# workbook, mirroring the existing Add/Delete sheet pairs already present.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. DeleteLocations sheet (pattern copied from the existing DeleteJobTitles
#    two-column Add/Validation sheet so it inherits the same styles/borders).
# ---------------------------------------------------------------------------
$refDelete = $wb.Worksheets.Item("DeleteJobTitles")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$refDelete.Copy($null, $lastSheet)
$delLoc = $wb.Worksheets.Item($wb.Worksheets.Count)
$delLoc.Name = "DeleteLocations"

$delLoc.Range("A1").Value = "locations"
$delLoc.Range("B1").Value = "Validation"
$delLoc.Range("A2").Value = "BENGALURU"
$delLoc.Range("B2").Value = $false
$delLoc.Range("A3").Value = "CHENNAI"
$delLoc.Range("B3").Value = $true

$delLoc.Columns.Item(1).ColumnWidth = 11.7109375
$delLoc.Range("A1:B3").Select()
$delLoc.Range("C1").Select()

# ---------------------------------------------------------------------------
# 2. AddLocation sheet (9-column location detail form).
# ---------------------------------------------------------------------------
$refAdd = $wb.Worksheets.Item("AddJobTitles")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$addLoc = $wb.Worksheets.Add($null, $lastSheet2)
$addLoc.Name = "AddLocation"

# -- header row (style1: bold/fill/full border) applied to A1:I1 first --
$refAdd.Range("A1").Copy()
$addLoc.Range("A1:I1").PasteSpecial(-4122)

# -- body rows (style2: plain + full border) applied to A2:I5 --
$refAdd.Range("A2").Copy()
$addLoc.Range("A2:I5").PasteSpecial(-4122)

# -- style5: centered alignment (used by the Country column body cells) --
$addLoc.Range("A2").Copy()
$addLoc.Range("K1").PasteSpecial(-4122)
$addLoc.Range("K1").HorizontalAlignment = -4108
$addLoc.Range("K1").VerticalAlignment = -4108
$addLoc.Range("K1").Copy()
$addLoc.Range("B2:B4").PasteSpecial(-4122)
$addLoc.Range("K1").Clear()

# -- style7 then style6: header border variants (left/right/top, then
#    left/right only) built by trimming the full header border --
$refAdd.Range("A1").Copy()
$addLoc.Range("K1").PasteSpecial(-4122)
$addLoc.Range("K1").Borders.Item(9).LineStyle = -4142
$addLoc.Range("K1").Copy()
$addLoc.Range("C1:H1").PasteSpecial(-4122)

$addLoc.Range("K1").Borders.Item(8).LineStyle = -4142
$addLoc.Range("K1").Copy()
$addLoc.Range("I1").PasteSpecial(-4122)
$addLoc.Range("K1").Clear()

# -- values --
$addLoc.Range("A1").Value = "Locations"
$addLoc.Range("B1").Value = "Country"
$addLoc.Range("C1").Value = "State"
$addLoc.Range("D1").Value = "City"
$addLoc.Range("E1").Value = "Address"
$addLoc.Range("F1").Value = "PinCode"
$addLoc.Range("G1").Value = "Phone"
$addLoc.Range("H1").Value = "Fax"
$addLoc.Range("I1").Value = "Notes"

$addLoc.Range("A2").Value = "CHENNAI"
$addLoc.Range("B2").Value = "India"
$addLoc.Range("C2").Value = "TAMILNADU"
$addLoc.Range("D2").Value = "CHENNAI"
$addLoc.Range("E2").Value = "chennai"
$addLoc.Range("F2").Value = 570985
$addLoc.Range("G2").Value = 999999999
$addLoc.Range("H2").Value = "NA"
$addLoc.Range("I2").Value = "Test Company"

$addLoc.Range("A3").Value = "Los Angles"
$addLoc.Range("B3").Value = "Algeria"

$addLoc.Range("B4").Value = "India"

$addLoc.Range("A5").Value = "MUMBAI"

# -- column widths --
$addLoc.Columns.Item(1).ColumnWidth = 11.7109375
$addLoc.Columns.Item(2).ColumnWidth = 12.85546875
$addLoc.Columns.Item(3).ColumnWidth = 11.85546875
$addLoc.Columns.Item(7).ColumnWidth = 10
$addLoc.Columns.Item(9).ColumnWidth = 13.5703125

$addLoc.Range("A2").Select()
$addLoc.Activate()

# ---------------------------------------------------------------------------
# 3. Misc selection tweaks on existing sheets captured in the diff.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("DeleteJobCategories").Range("A1:B3").Select()
$wb.Worksheets.Item("DeleteJobCategories").Range("B3").Activate()

$addLoc.Activate()

Write-Host "done"
